$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Publication date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact / No display for ContactDetail" row (row 11),
# then turn the remaining row 10 into "Jurisdiction / United States of America"
$ws.Rows.Item(11).Delete()
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Elements sheet: root Extension row Short/Definition now describe "Path"
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "Path"
$ws2.Range("L2").Value = "Path to a FHIR element"
